# Update the crypto price/volume table with refreshed data.
# Mirrors the scheduled "Updated cryptos list ... with GitHub Actions" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting
# (values like "1.000" or "27.676.10" are not real numbers) by forcing
# the cells to Text format before writing the new values.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.669.98"
$ws.Range("E2").Value = "  -0.77%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.898.22"
$ws.Range("E3").Value = "  -0.34%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.17%  "

# Row 5 - BNB
$ws.Range("D5").Value = "310.59"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.08%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +5.04%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3802"
$ws.Range("E8").Value = "  -0.55%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.07239"
$ws.Range("E9").Value = "  -1.40%  "

# Row 10 - Solana (was Polygon)
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "21.09"
$ws.Range("E10").Value = "  +1.00%  "

# Row 11 - Polygon (was Solana)
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "0.9024"
$ws.Range("E11").Value = "  -0.97%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.888.13"
$ws.Range("E12").Value = "  -1.41%  "

# Row 13 - TRON
$ws.Range("D13").Value = "0.07636"
$ws.Range("E13").Value = "  -0.04%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "5.437"
$ws.Range("E14").Value = "  -1.06%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "91.67"
$ws.Range("E15").Value = "  +0.26%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  -0.11%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.000008658"
$ws.Range("E17").Value = "  -0.72%  "

# Row 18 - Avalanche
$ws.Range("E18").Value = "  -1.36%  "

# Row 19 - Dai
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.18%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "27.688.10"
$ws.Range("E20").Value = "  -0.83%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "5.150"
$ws.Range("E21").Value = "  +0.33%  "

# Row 22 - WrappedliquidstakedEther2.0
$ws.Range("D22").Value = "2.121.57"
$ws.Range("E22").Value = "  -0.69%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -0.09%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "6.603"
$ws.Range("E24").Value = "  -0.49%  "

# Row 25 - Monero
$ws.Range("D25").Value = "153.46"
$ws.Range("E25").Value = "  -0.73%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -0.50%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "18.27"
$ws.Range("E27").Value = "  -0.73%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "2.199"
$ws.Range("E28").Value = "  -1.53%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "114.26"
$ws.Range("E29").Value = "  -1.00%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "4.833"
$ws.Range("E30").Value = "  -2.04%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "4.809"
$ws.Range("E31").Value = "  +3.38%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "0.09149"
$ws.Range("E32").Value = "  +1.73%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.05270"
$ws.Range("E33").Value = "  -0.24%  "

# Row 34 - HuobiToken
$ws.Range("D34").Value = "3.123"
$ws.Range("E34").Value = "  -2.46%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "1.221"
$ws.Range("E35").Value = "  -1.50%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.7708"
$ws.Range("E36").Value = "  -0.44%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.02087"
$ws.Range("E37").Value = "  +1.07%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "2.570"
$ws.Range("E38").Value = "  +0.12%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "3.075"
$ws.Range("E39").Value = "  +2.02%  "

# Row 40 - TheSandbox (was TrustWalletToken)
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.5573"
$ws.Range("E40").Value = "  +0.79%  "

# Row 41 - TrustWalletToken (was TheSandbox)
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.091"
$ws.Range("E41").Value = "  -0.92%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "6.713"
$ws.Range("E42").Value = "  -3.91%  "

# Row 43 - Quant
$ws.Range("D43").Value = "117.26"
$ws.Range("E43").Value = "  +5.02%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "8.712"
$ws.Range("E44").Value = "  +1.80%  "

# Row 45 - Algorand
$ws.Range("E45").Value = "  -0.95%  "

# Row 46 - Decentraland
$ws.Range("D46").Value = "0.4806"
$ws.Range("E46").Value = "  +0.11%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "10.36"
$ws.Range("E47").Value = "  -2.56%  "

# Row 48 - PaxDollar
$ws.Range("D48").Value = "1.0000"
$ws.Range("E48").Value = "  -0.09%  "

# Row 49 - NEARProtocol
$ws.Range("E49").Value = "  -2.97%  "

# Row 50 - Aave
$ws.Range("E50").Value = "  -1.92%  "

# Row 51 - Elrond
$ws.Range("D51").Value = "37.10"
$ws.Range("E51").Value = "  +0.15%  "
